# Refractor FRAME to DTI
# Update the byline placeholders on the title slide (slide 1):
#   - "October 2021"  -> "Wireless Specialist (PACNW)" plus a trailing blank line
#   - "JR & RW"        -> "Randy Wu, Principal ESE"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shapes = $s.Shapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $shp = $shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }

    $tr = $shp.TextFrame.TextRange
    $txt = $tr.Text

    if ($txt -eq "October 2021") {
        $tr.Text = "Wireless Specialist (PACNW)`r"
    }
    elseif ($txt -eq "JR & RW") {
        $tr.Text = "Randy Wu, Principal ESE"
    }
}
